$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every cell we are about to touch so that
# numeric-looking / date-looking strings (e.g. "244.01", "2-1-2023")
# are preserved verbatim as text instead of being parsed into
# numbers or dates by Excel. NumberFormat is applied to a
# contiguous span (single-area range) BEFORE the values are
# written, since multi-area (union) ranges do not reliably
# propagate the format to every cell in this host.

# Row 2
$ws.Range("D2:G2").NumberFormat = "@"
$ws.Range("D2").Value = "244.01"
$ws.Range("E2").Value = "-0.88%"
$ws.Range("F2").Value = "2-1-2023"
$ws.Range("G2").Value = "0"

# Row 3
$ws.Range("D3:G3").NumberFormat = "@"
$ws.Range("D3").Value = "27.31"
$ws.Range("E3").Value = "5.44%"
$ws.Range("F3").Value = "2-1-2023"
$ws.Range("G3").Value = "0"

# Row 4
$ws.Range("D4:G4").NumberFormat = "@"
$ws.Range("D4").Value = "5.097"
$ws.Range("E4").Value = "-1.49%"
$ws.Range("F4").Value = "2-1-2023"
$ws.Range("G4").Value = "0"

# Row 5
$ws.Range("D5:G5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05664"
$ws.Range("E5").Value = "1.49%"
$ws.Range("F5").Value = "2-1-2023"
$ws.Range("G5").Value = "0"

# Row 6
$ws.Range("D6:G6").NumberFormat = "@"
$ws.Range("D6").Value = "6.536"
$ws.Range("E6").Value = "0.88%"
$ws.Range("F6").Value = "2-1-2023"
$ws.Range("G6").Value = "0"

# Row 7
$ws.Range("D7:G7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8210"
$ws.Range("E7").Value = "0.94%"
$ws.Range("F7").Value = "2-1-2023"
$ws.Range("G7").Value = "0"

# Row 8
$ws.Range("D8:G8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8520"
$ws.Range("E8").Value = "1.13%"
$ws.Range("F8").Value = "2-1-2023"
$ws.Range("G8").Value = "0"

# Row 9
$ws.Range("D9:G9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1334"
$ws.Range("E9").Value = "0.41%"
$ws.Range("F9").Value = "2-1-2023"
$ws.Range("G9").Value = "0"

# Row 10
$ws.Range("B10:G10").NumberFormat = "@"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.06942"
$ws.Range("E10").Value = "0.03%"
$ws.Range("F10").Value = "2-1-2023"
$ws.Range("G10").Value = "0"

# Row 11
$ws.Range("B11:G11").NumberFormat = "@"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.02870"
$ws.Range("E11").Value = "1.22%"
$ws.Range("F11").Value = "2-1-2023"
$ws.Range("G11").Value = "0"

# Row 12
$ws.Range("B12:G12").NumberFormat = "@"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "0.09386"
$ws.Range("E12").Value = "0.06%"
$ws.Range("F12").Value = "2-1-2023"
$ws.Range("G12").Value = "0"

# Row 13
$ws.Range("B13:G13").NumberFormat = "@"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "0.001514"
$ws.Range("E13").Value = "-0.55%"
$ws.Range("F13").Value = "2-1-2023"
$ws.Range("G13").Value = "0"

# Row 14
$ws.Range("B14:G14").NumberFormat = "@"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").Value = "0.04119"
$ws.Range("E14").Value = "-10.97%"
$ws.Range("F14").Value = "2-1-2023"
$ws.Range("G14").Value = "0"

# Row 15
$ws.Range("B15:G15").NumberFormat = "@"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0006013"
$ws.Range("E15").Value = "1.04%"
$ws.Range("F15").Value = "2-1-2023"
$ws.Range("G15").Value = "0"

# Row 16
$ws.Range("B16:G16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006219"
$ws.Range("E16").Value = "1.35%"
$ws.Range("F16").Value = "2-1-2023"
$ws.Range("G16").Value = "0"

# Row 17
$ws.Range("B17:G17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.514"
$ws.Range("E17").Value = "-2.61%"
$ws.Range("F17").Value = "2-1-2023"
$ws.Range("G17").Value = "0"

# Row 18
$ws.Range("B18:G18").NumberFormat = "@"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "3.018"
$ws.Range("E18").Value = "-0.32%"
$ws.Range("F18").Value = "2-1-2023"
$ws.Range("G18").Value = "0"

# Row 19
$ws.Range("B19:G19").NumberFormat = "@"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.202"
$ws.Range("E19").Value = "0.87%"
$ws.Range("F19").Value = "2-1-2023"
$ws.Range("G19").Value = "0"

# Row 20
$ws.Range("B20:G20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3149"
$ws.Range("E20").Value = "1.21%"
$ws.Range("F20").Value = "2-1-2023"
$ws.Range("G20").Value = "0"

# Row 21
$ws.Range("D21:G21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03288"
$ws.Range("E21").Value = "3.89%"
$ws.Range("F21").Value = "2-1-2023"
$ws.Range("G21").Value = "0"

# Row 22
$ws.Range("E22:G22").NumberFormat = "@"
$ws.Range("E22").Value = "2.41%"
$ws.Range("F22").Value = "2-1-2023"
$ws.Range("G22").Value = "0"

# Row 23
$ws.Range("D23:G23").NumberFormat = "@"
$ws.Range("D23").Value = "3.575"
$ws.Range("E23").Value = "-4.37%"
$ws.Range("F23").Value = "2-1-2023"
$ws.Range("G23").Value = "0"

# Row 24
$ws.Range("D24:G24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1375"
$ws.Range("E24").Value = "0.04%"
$ws.Range("F24").Value = "2-1-2023"
$ws.Range("G24").Value = "0"

# Row 25
$ws.Range("D25:G25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001210"
$ws.Range("E25").Value = "-2.84%"
$ws.Range("F25").Value = "2-1-2023"
$ws.Range("G25").Value = "0"

# Row 26
$ws.Range("D26:G26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004464"
$ws.Range("E26").Value = "-1.50%"
$ws.Range("F26").Value = "2-1-2023"
$ws.Range("G26").Value = "0"

# Row 27
$ws.Range("E27:G27").NumberFormat = "@"
$ws.Range("E27").Value = "22.88%"
$ws.Range("F27").Value = "2-1-2023"
$ws.Range("G27").Value = "0"

# Row 28
$ws.Range("E28:G28").NumberFormat = "@"
$ws.Range("E28").Value = "0.36%"
$ws.Range("F28").Value = "2-1-2023"
$ws.Range("G28").Value = "0"

# Row 29
$ws.Range("F29:G29").NumberFormat = "@"
$ws.Range("F29").Value = "2-1-2023"
$ws.Range("G29").Value = "0"

# Row 30
$ws.Range("F30:G30").NumberFormat = "@"
$ws.Range("F30").Value = "2-1-2023"
$ws.Range("G30").Value = "0"

# Row 31
$ws.Range("F31:G31").NumberFormat = "@"
$ws.Range("F31").Value = "2-1-2023"
$ws.Range("G31").Value = "0"

# Row 32
$ws.Range("F32:G32").NumberFormat = "@"
$ws.Range("F32").Value = "2-1-2023"
$ws.Range("G32").Value = "0"

# Row 33
$ws.Range("F33:G33").NumberFormat = "@"
$ws.Range("F33").Value = "2-1-2023"
$ws.Range("G33").Value = "0"

# Row 34
$ws.Range("F34:G34").NumberFormat = "@"
$ws.Range("F34").Value = "2-1-2023"
$ws.Range("G34").Value = "0"

# Row 35
$ws.Range("F35:G35").NumberFormat = "@"
$ws.Range("F35").Value = "2-1-2023"
$ws.Range("G35").Value = "0"

# Row 36
$ws.Range("F36:G36").NumberFormat = "@"
$ws.Range("F36").Value = "2-1-2023"
$ws.Range("G36").Value = "0"

# Row 37
$ws.Range("F37:G37").NumberFormat = "@"
$ws.Range("F37").Value = "2-1-2023"
$ws.Range("G37").Value = "0"

# Row 38
$ws.Range("F38:G38").NumberFormat = "@"
$ws.Range("F38").Value = "2-1-2023"
$ws.Range("G38").Value = "0"

# Row 39
$ws.Range("F39:G39").NumberFormat = "@"
$ws.Range("F39").Value = "2-1-2023"
$ws.Range("G39").Value = "0"

# Row 40
$ws.Range("D40:G40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03708"
$ws.Range("E40").Value = "1.66%"
$ws.Range("F40").Value = "2-1-2023"
$ws.Range("G40").Value = "0"

# Row 41
$ws.Range("D41:G41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005950"
$ws.Range("E41").Value = "-3.28%"
$ws.Range("F41").Value = "2-1-2023"
$ws.Range("G41").Value = "0"

# Row 42
$ws.Range("D42:G42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1056"
$ws.Range("E42").Value = "0.44%"
$ws.Range("F42").Value = "2-1-2023"
$ws.Range("G42").Value = "0"

# Row 43
$ws.Range("D43:G43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002313"
$ws.Range("E43").Value = "-7.51%"
$ws.Range("F43").Value = "2-1-2023"
$ws.Range("G43").Value = "0"

# Row 44
$ws.Range("D44:G44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009055"
$ws.Range("E44").Value = "13.12%"
$ws.Range("F44").Value = "2-1-2023"
$ws.Range("G44").Value = "0"

# Row 45
$ws.Range("E45:G45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.58%"
$ws.Range("F45").Value = "2-1-2023"
$ws.Range("G45").Value = "0"

# Row 46
$ws.Range("D46:G46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("F46").Value = "2-1-2023"
$ws.Range("G46").Value = "0"

# Row 47
$ws.Range("D47:G47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1191"
$ws.Range("E47").Value = "-17.90%"
$ws.Range("F47").Value = "2-1-2023"
$ws.Range("G47").Value = "0"

# Row 48
$ws.Range("D48:G48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002523"
$ws.Range("E48").Value = "4.83%"
$ws.Range("F48").Value = "2-1-2023"
$ws.Range("G48").Value = "0"

# Row 49
$ws.Range("F49:G49").NumberFormat = "@"
$ws.Range("F49").Value = "2-1-2023"
$ws.Range("G49").Value = "0"

# Row 50
$ws.Range("F50:G50").NumberFormat = "@"
$ws.Range("F50").Value = "2-1-2023"
$ws.Range("G50").Value = "0"

# Row 51
$ws.Range("F51:G51").NumberFormat = "@"
$ws.Range("F51").Value = "2-1-2023"
$ws.Range("G51").Value = "0"

